{"js": "const replacements = [\n  [\"1.174\", \"0.019\"],\n  [\"0.006\", \"0.009\"],\n  [\"3.457\", \"4.831\"],\n  [\"0.027\", \"0.008\"],\n  [\"1.772\", \"0.015\"],\n  [\"0.010\", \"0.007\"],\n  [\"1.739\", \"1.238\"],\n  [\"0.054\", \"0.000\"],\n  [\"0.320\", \"0.039\"],\n  [\"0.199\", \"0.006\"],\n  [\"0.001\", \"0.003\"],\n  [\"0.293\", \"0.776\"],\n  [\"178.158\", \"2.072\"],\n  [\"0.982\", \"0.981\"],\n  [\"181.357\", \"2.112\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(`Expected exactly one match for \"${oldText}\", found ${results.items.length}`);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"1.174\";   New = \"0.019\" },\n    @{ Old = \"0.006\";   New = \"0.009\" },\n    @{ Old = \"3.457\";   New = \"4.831\" },\n    @{ Old = \"0.027\";   New = \"0.008\" },\n    @{ Old = \"1.772\";   New = \"0.015\" },\n    @{ Old = \"0.010\";   New = \"0.007\" },\n    @{ Old = \"1.739\";   New = \"1.238\" },\n    @{ Old = \"0.054\";   New = \"0.000\" },\n    @{ Old = \"0.320\";   New = \"0.039\" },\n    @{ Old = \"0.199\";   New = \"0.006\" },\n    @{ Old = \"0.001\";   New = \"0.003\" },\n    @{ Old = \"0.293\";   New = \"0.776\" },\n    @{ Old = \"178.158\"; New = \"2.072\" },\n    @{ Old = \"0.982\";   New = \"0.981\" },\n    @{ Old = \"181.357\"; New = \"2.112\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap(wdFindContinue=1), Format, ReplaceWith,\n    # Replace(wdReplaceAll=2)\n    $find.Execute(\n        $r.Old, $true, $true, $false, $false, $false, $true, 1, $false,\n        $r.New, 2\n    )\n}\n"}
